$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 3).Value = 6.128577557148986
$ws.Cells.Item(18, 3).Value = 6.1766
$ws.Cells.Item(19, 3).Value = 6.2847
$ws.Cells.Item(20, 3).Value = 6.3757
$ws.Cells.Item(22, 3).Value = 6.7736909842173
$ws.Cells.Item(23, 3).Value = 6.7736909842173
$ws.Cells.Item(24, 3).Value = 6.7736909842173
$ws.Cells.Item(25, 3).Value = 6.873276000000001
$ws.Cells.Item(26, 3).Value = 7.0058

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(16, 3).Value = 6.128577557148986
$ws.Cells.Item(17, 3).Value = 6.309950000000001
$ws.Cells.Item(19, 3).Value = 6.7736909842173
$ws.Cells.Item(20, 3).Value = 6.7736909842173
$ws.Cells.Item(21, 3).Value = 6.801851
$ws.Cells.Item(22, 3).Value = 7.0616

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(17, 3).Value = 6.128577557148986
$ws.Cells.Item(18, 3).Value = 6.2136
$ws.Cells.Item(19, 3).Value = 6.3504
$ws.Cells.Item(22, 3).Value = 6.7736909842173
$ws.Cells.Item(23, 3).Value = 6.886656
$ws.Cells.Item(24, 3).Value = 6.95474
$ws.Cells.Item(25, 3).Value = 7.046490243855914

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(19, 3).Value = 6.13815
$ws.Cells.Item(20, 3).Value = 6.2285
$ws.Cells.Item(21, 3).Value = 6.3211
$ws.Cells.Item(22, 3).Value = 6.39
$ws.Cells.Item(25, 3).Value = 6.7736909842173
$ws.Cells.Item(26, 3).Value = 6.86982
$ws.Cells.Item(27, 3).Value = 6.998527437772841

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(20, 3).Value = 6.1447
$ws.Cells.Item(21, 3).Value = 6.2558
$ws.Cells.Item(22, 3).Value = 6.3599
$ws.Cells.Item(27, 3).Value = 6.822156000000001
$ws.Cells.Item(28, 3).Value = 6.97771
$ws.Cells.Item(29, 3).Value = 7.053795

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(20, 3).Value = 6.128577557148986
$ws.Cells.Item(21, 3).Value = 6.1955
$ws.Cells.Item(22, 3).Value = 6.3
$ws.Cells.Item(23, 3).Value = 6.37505
$ws.Cells.Item(27, 3).Value = 6.7736909842173
$ws.Cells.Item(28, 3).Value = 6.88826
$ws.Cells.Item(29, 3).Value = 6.990284022413692
$ws.Cells.Item(30, 3).Value = 7.073555040636824

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 3).Value = 6.128577557148986
$ws.Cells.Item(23, 3).Value = 6.24655
$ws.Cells.Item(24, 3).Value = 6.358827391277695
$ws.Cells.Item(30, 3).Value = 6.7736909842173
$ws.Cells.Item(31, 3).Value = 6.824652172142002
$ws.Cells.Item(32, 3).Value = 6.888862689544663
$ws.Cells.Item(33, 3).Value = 6.951375078310891
$ws.Cells.Item(34, 3).Value = 7.012764308233089
$ws.Cells.Item(35, 3).Value = 7.073605349103663

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(21, 3).Value = 6.128577557148986
$ws.Cells.Item(22, 3).Value = 6.21585
$ws.Cells.Item(23, 3).Value = 6.3263
$ws.Cells.Item(24, 3).Value = 6.39935
$ws.Cells.Item(28, 3).Value = 6.7736909842173
$ws.Cells.Item(29, 3).Value = 6.939438328537795
$ws.Cells.Item(30, 3).Value = 7.039263628647083

$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(18, 3).Value = 6.128577557148986
$ws.Cells.Item(19, 3).Value = 6.21565
$ws.Cells.Item(20, 3).Value = 6.3385
$ws.Cells.Item(23, 3).Value = 6.7736909842173
$ws.Cells.Item(24, 3).Value = 6.7736909842173
$ws.Cells.Item(25, 3).Value = 6.876046
$ws.Cells.Item(26, 3).Value = 6.962471
$ws.Cells.Item(27, 3).Value = 7.07656

$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(16, 3).Value = 6.128577557148986
$ws.Cells.Item(17, 3).Value = 6.1795
$ws.Cells.Item(18, 3).Value = 6.2972
$ws.Cells.Item(19, 3).Value = 6.4121
$ws.Cells.Item(21, 3).Value = 6.7736909842173
$ws.Cells.Item(22, 3).Value = 6.7736909842173
$ws.Cells.Item(23, 3).Value = 6.830762
$ws.Cells.Item(24, 3).Value = 6.957697
$ws.Cells.Item(25, 3).Value = 7.047104

$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(16, 3).Value = 6.128577557148986
$ws.Cells.Item(17, 3).Value = 6.2662
$ws.Cells.Item(18, 3).Value = 6.3838
$ws.Cells.Item(20, 3).Value = 6.7736909842173
$ws.Cells.Item(21, 3).Value = 6.7736909842173
$ws.Cells.Item(22, 3).Value = 6.84633
$ws.Cells.Item(23, 3).Value = 7.003941999999999
$ws.Cells.Item(24, 3).Value = 7.066104

$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(16, 3).Value = 6.153700000000001
$ws.Cells.Item(17, 3).Value = 6.2892
$ws.Cells.Item(20, 3).Value = 6.7736909842173
$ws.Cells.Item(21, 3).Value = 6.885344
$ws.Cells.Item(22, 3).Value = 7.024208000000001

$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(15, 3).Value = 6.128577557148986
$ws.Cells.Item(16, 3).Value = 6.18275
$ws.Cells.Item(17, 3).Value = 6.389
$ws.Cells.Item(18, 3).Value = 6.7736909842173
$ws.Cells.Item(19, 3).Value = 6.7736909842173
$ws.Cells.Item(20, 3).Value = 6.809672
$ws.Cells.Item(21, 3).Value = 6.956778
